$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: paragraph "People will install an monster that passively
# grows off of their sms usage" -> "Will people install an app that
# passively grows off of their sms usage?"
# (the "sms" run keeps its spell-check proofErr wrapper, so only
#  replace the text around it)
# -----------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$rng1 = $p6.Range
$rng1.Find.Execute(
    "People will install an monster that passively grows off of their ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Will people install an app that passively grows off of their ",
    2)

$p6b = $d.Paragraphs.Item(6)
$rng1b = $p6b.Range
$rng1b.Find.Execute(" usage", $true, $false, $false, $false, $false, $true, 1, $false, " usage?", 2)

# -----------------------------------------------------------------
# Change 2: split the "4) People would ..." paragraph's single run
# into "4) " and "People would ..." runs, and add a _GoBack bookmark
# right after the paragraph text (before the paragraph mark).
# -----------------------------------------------------------------
$p19 = $d.Paragraphs.Item(19)
$pStart = $p19.Range.Start

# Force a run split right after "4) " using a transient bookmark.
$splitPoint = $d.Range($pStart + 3, $pStart + 3)
$d.Bookmarks.Add("TempSplit", $splitPoint)
$d.Bookmarks.Item("TempSplit").Delete()

# Add the _GoBack bookmark at the very end of the paragraph text
# (immediately before the paragraph mark). A collapsed range exactly
# at that boundary is placed incorrectly, so insert a throw-away
# character after the target point first, anchor the bookmark before
# it, then remove the throw-away character.
$p19b = $d.Paragraphs.Item(19)
$rngEnd = $p19b.Range.End
$insertPos = $d.Range($rngEnd - 1, $rngEnd - 1)
$insertPos.InsertAfter("X")
$target = $d.Range($rngEnd - 1, $rngEnd - 1)
$d.Bookmarks.Add("_GoBack", $target)
$xRange = $d.Range($rngEnd - 1, $rngEnd)
$xRange.Delete()
